$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2 updates
$ws.Range("J2").Value = 3.85
$ws.Range("M2").Value = 1.05
$ws.Range("N2").Value = 3.7
$ws.Range("P2").Value = 2.1
$ws.Range("Q2").Value = 1.64
$ws.Range("R2").Value = 1.43

# Row 3 updates
$ws.Range("J3").Value = 4.9
$ws.Range("L3").Value = 1.29
$ws.Range("T3").Value = 1.85
$ws.Range("U3").Value = 1.96
$ws.Range("Y3").Value = 12
$ws.Range("AG3").Value = 38
